$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")

# Row 28
$ws.Range("H28").Value = 1247.1666
$ws.Range("I28").Value = 1188.3846
$ws.Range("J28").Value = 1400
$ws.Range("K28").Value = 1188.3846
$ws.Range("L28").Value = 1400
$ws.Range("M28").Value = -703.3846000000001
$ws.Range("N28").Value = -2370

# Row 51
$ws.Range("H51").Value = 1854449.4
$ws.Range("I51").Value = 5556405.5
$ws.Range("J51").Value = 3471.25
$ws.Range("K51").Value = 5556405.5
$ws.Range("L51").Value = 3471.25
$ws.Range("M51").Value = -5555921.5
$ws.Range("N51").Value = -4439.25

# Row 61
$ws.Range("H61").Value = 568
$ws.Range("I61").Value = 568
$ws.Range("K61").Value = 1704
$ws.Range("M61").Value = -1532

# Row 62
$ws.Range("H62").Value = 1367.1111
$ws.Range("I62").Value = 1367.1111
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 1367.1111
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents() | Out-Null
$ws.Range("N62").Value = -743.1111000000001

# Row 65
$ws.Range("H65").Value = 1367.1111
$ws.Range("I65").Value = 1367.1111
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 6835.5555
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -3715.5555
$ws.Range("N65").Value = 0

# Row 94
$ws.Range("H94").Value = 8000
$ws.Range("I94").Value = 8000
$ws.Range("K94").Value = 8000
$ws.Range("M94").Value = -7549

# Row 112
$ws.Range("H112").Value = 2120.5
$ws.Range("I112").Value = 416.66666
$ws.Range("J112").Value = 2285.3872
$ws.Range("K112").Value = 1249.99998
$ws.Range("L112").Value = 6856.1616
$ws.Range("M112").Value = -141.9999800000001
$ws.Range("N112").Value = -9072.161599999999

# Row 131
$ws.Range("H131").Value = 501797.5
$ws.Range("I131").Value = 1000095
$ws.Range("J131").Value = 3500
$ws.Range("K131").Value = 3000285
$ws.Range("L131").Value = 10500
$ws.Range("M131").Value = -2995245
$ws.Range("N131").Value = -20580


# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")

# Row 32
$ws.Range("H32").Value = 8489.482
$ws.Range("I32").Value = 4241.771
$ws.Range("J32").Value = 28878.5
$ws.Range("K32").Value = 4241.771
$ws.Range("L32").Value = 28878.5
$ws.Range("M32").Value = -3954.771
$ws.Range("N32").Value = -29452.5

# Row 61
$ws.Range("H61").Value = 1425.5
$ws.Range("I61").Value = 1186.0834
$ws.Range("J61").Value = 1712.8
$ws.Range("K61").Value = 1186.0834
$ws.Range("L61").Value = 1712.8
$ws.Range("M61").Value = -974.0834
$ws.Range("N61").Value = -2136.8

# Row 74
$ws.Range("H74").Value = 83335256
$ws.Range("I74").Value = 125001500
$ws.Range("J74").Value = 2750
$ws.Range("K74").Value = 125001500
$ws.Range("L74").Value = 2750
$ws.Range("M74").Value = -125000626
$ws.Range("N74").Value = -4498

# Row 77
$ws.Range("H77").Value = 83335256
$ws.Range("I77").Value = 125001500
$ws.Range("J77").Value = 2750
$ws.Range("K77").Value = 625007500
$ws.Range("L77").Value = 13750
$ws.Range("M77").Value = -625003132
$ws.Range("N77").Value = -22486

# Row 109
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").ClearContents() | Out-Null
$ws.Range("N109").Value = 0

# Row 122
$ws.Range("H122").Value = 2477.261
$ws.Range("I122").Value = 2630
$ws.Range("J122").Value = 2359.7693
$ws.Range("K122").Value = 7890
$ws.Range("L122").Value = 7079.3079
$ws.Range("M122").Value = -5440
$ws.Range("N122").Value = -11979.3079

# Row 135
$ws.Range("H135").Value = 25725
$ws.Range("J135").Value = 25725
$ws.Range("L135").Value = 25725
$ws.Range("N135").Value = -35865

# Row 136
$ws.Range("H136").Value = 1425.5
$ws.Range("I136").Value = 1186.0834
$ws.Range("J136").Value = 1712.8
$ws.Range("K136").Value = 3558.2502
$ws.Range("L136").Value = 5138.4
$ws.Range("M136").Value = -1008.2502
$ws.Range("N136").Value = -10238.4


# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")

# Row 98
$ws.Range("H98").Value = 57250
$ws.Range("J98").Value = 57250
$ws.Range("L98").Value = 57250
$ws.Range("N98").Value = -63240

# Row 102
$ws.Range("H102").Value = 16728
$ws.Range("I102").Value = 16728
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 16728
$ws.Range("L102").Value = 0
$ws.Range("M102").ClearContents() | Out-Null
$ws.Range("N102").Value = -13483

# Row 105
$ws.Range("H105").Value = 4547796.5
$ws.Range("I105").Value = 11366136
$ws.Range("J105").Value = 2237
$ws.Range("K105").Value = 11366136
$ws.Range("L105").Value = 2237
$ws.Range("M105").Value = -11364389
$ws.Range("N105").Value = -5731

# Row 134
$ws.Range("H134").Value = 6411243.5
$ws.Range("I134").Value = 7143805.5
$ws.Range("J134").Value = 1325
$ws.Range("K134").Value = 21431416.5
$ws.Range("L134").Value = 3975
$ws.Range("M134").Value = -21428881.5
$ws.Range("N134").Value = -9045


# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")

# Row 31
$ws.Range("H31").Value = 14708701
$ws.Range("I31").Value = 20835230
$ws.Range("J31").Value = 5029.3
$ws.Range("K31").Value = 20835230
$ws.Range("L31").Value = 5029.3
$ws.Range("M31").Value = -20834935
$ws.Range("N31").Value = -5619.3

# Row 34
$ws.Range("H34").Value = 14708701
$ws.Range("I34").Value = 20835230
$ws.Range("J34").Value = 5029.3
$ws.Range("K34").Value = 20835230
$ws.Range("L34").Value = 5029.3
$ws.Range("M34").Value = -20835028
$ws.Range("N34").Value = -5433.3

# Row 94
$ws.Range("H94").Value = 2642.5
$ws.Range("I94").Value = 3403
$ws.Range("J94").Value = 2338.3
$ws.Range("K94").Value = 3403
$ws.Range("L94").Value = 2338.3
$ws.Range("M94").Value = -2952
$ws.Range("N94").Value = -3240.3

# Row 111
$ws.Range("H111").Value = 55351
$ws.Range("J111").Value = 55351
$ws.Range("L111").Value = 55351
$ws.Range("N111").Value = -63531

# Row 141
$ws.Range("H141").Value = 164333.33
$ws.Range("J141").Value = 164333.33
$ws.Range("L141").Value = 164333.33
$ws.Range("N141").Value = -174693.33


# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")

# Row 3
$ws.Range("H3").Value = 3491.2
$ws.Range("I3").Value = 3235
$ws.Range("K3").Value = 9705
$ws.Range("M3").Value = -9593

# Row 95
$ws.Range("H95").Value = 8333.333000000001
$ws.Range("J95").Value = 8333.333000000001
$ws.Range("L95").Value = 24999.999
$ws.Range("N95").Value = -29117.999

# Row 102
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").ClearContents() | Out-Null
$ws.Range("N102").Value = 0


# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")

# Row 102
$ws.Range("H102").Value = 36357.344
$ws.Range("I102").Value = 48873.76
$ws.Range("J102").Value = 3501.75
$ws.Range("K102").Value = 48873.76
$ws.Range("L102").Value = 3501.75
$ws.Range("M102").Value = -47251.76
$ws.Range("N102").Value = -6745.75

# Row 123
$ws.Range("H123").Value = 11016.889
$ws.Range("J123").Value = 11016.889
$ws.Range("L123").Value = 11016.889
$ws.Range("N123").Value = -15916.889


# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")

# Row 133
$ws.Range("H133").Value = 23775
$ws.Range("J133").Value = 23775
$ws.Range("L133").Value = 23775
$ws.Range("N133").Value = -28835


# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")

# Row 136
$ws.Range("H136").Value = 5113.143
$ws.Range("I136").Value = 1168.3158
$ws.Range("J136").Value = 13441.111
$ws.Range("K136").Value = 3504.9474
$ws.Range("L136").Value = 40323.333
$ws.Range("M136").Value = -954.9474
$ws.Range("N136").Value = -45423.333

